$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "307.16"
Set-TextValue "E2" "0.94%"
Set-TextValue "D3" "38.41"
Set-TextValue "E3" "7.49%"
Set-TextValue "D4" "5.083"
Set-TextValue "E4" "0.93%"
Set-TextValue "E5" "1.29%"
Set-TextValue "D6" "1.957"
Set-TextValue "E6" "5.67%"
Set-TextValue "B7" "KuCoinToken"
Set-TextValue "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D7" "7.951"
Set-TextValue "E7" "2.40%"
Set-TextValue "B8" "MXToken"
Set-TextValue "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.9276"
Set-TextValue "E8" "0.46%"
Set-TextValue "B9" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D9" "0.1459"
Set-TextValue "E9" "14.34%"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1933"
Set-TextValue "E10" "2.28%"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.09073"
Set-TextValue "E11" "0.23%"
Set-TextValue "B12" "BitrueCoin"
Set-TextValue "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03504"
Set-TextValue "E12" "2.71%"
Set-TextValue "B13" "BitMartToken"
Set-TextValue "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09800"
Set-TextValue "E13" "-0.66%"
Set-TextValue "B14" "BitForexToken"
Set-TextValue "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001412"
Set-TextValue "E14" "0.12%"
Set-TextValue "B15" "TigerCash"
Set-TextValue "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.006175"
Set-TextValue "E15" "-0.20%"
Set-TextValue "B16" "LEO"
Set-TextValue "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.723"
Set-TextValue "E16" "-3.32%"
Set-TextValue "B17" "GateToken"
Set-TextValue "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D17" "4.207"
Set-TextValue "E17" "1.95%"
Set-TextValue "D18" "3.416"
Set-TextValue "E18" "1.43%"
Set-TextValue "E19" "1.42%"
Set-TextValue "D20" "0.1334"
Set-TextValue "E20" "-0.90%"
Set-TextValue "D21" "4.806"
Set-TextValue "E21" "-0.14%"
Set-TextValue "D22" "0.2456"
Set-TextValue "E22" "6.34%"
Set-TextValue "D23" "0.04360"
Set-TextValue "E23" "-1.24%"
Set-TextValue "D24" "0.001228"
Set-TextValue "E24" "-0.57%"
Set-TextValue "E25" "-0.83%"
Set-TextValue "D27" "0.0001302"
Set-TextValue "E27" "-0.02%"
Set-TextValue "D39" "0.02085"
Set-TextValue "E39" "7.88%"
Set-TextValue "D40" "0.05060"
Set-TextValue "E40" "-1.15%"
Set-TextValue "D41" "0.007475"
Set-TextValue "E41" "-1.81%"
Set-TextValue "E42" "-0.24%"
Set-TextValue "D43" "0.1358"
Set-TextValue "E43" "0.63%"
Set-TextValue "D44" "0.002143"
Set-TextValue "E44" "-0.49%"
Set-TextValue "D45" "0.009200"
Set-TextValue "E45" "-6.93%"
Set-TextValue "D46" "0.00006189"
Set-TextValue "E46" "-0.10%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.00%"
Set-TextValue "D48" "0.002986"
Set-TextValue "D49" "0.001602"
Set-TextValue "E49" "-3.54%"
Set-TextValue "D50" "0.00002104"
Set-TextValue "E50" "0.00%"
Set-TextValue "D51" "0.0002003"
Set-TextValue "E51" "0.00%"
